# finish basic registry interaction
#
# 1) Tighten the table's left indent / left cell margins a couple of
#    twips (38 -> 36 dxa table indent; 33 -> 30 dxa cell margins).
# 2) Split "Finished basic log in ajax post, path redirection and
#    error dialog." into an explicit "1. " prefixed item and append a
#    new "2. ..." entry describing the finished registry work.
# 3) Mint the new (unused) ListLabel character styles 127-180 that a
#    save of this document carries along, mirroring the original
#    authoring tool's style bookkeeping.

$d = $word.ActiveDocument

# --- 1. Table spacing -------------------------------------------------
$t = $d.Tables(1)

# w:tblInd 38 -> 36 dxa  (1 dxa = 1/20 pt -> 1.9pt -> 1.8pt)
$t.Rows.LeftIndent = 1.8

# w:tblCellMar/left 33 -> 30 dxa (1.65pt -> 1.5pt)
$t.LeftPadding = 1.5

# w:tcMar/left 33 -> 30 dxa on every existing cell
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    for ($j = 1; $j -le $t.Columns.Count; $j++) {
        $t.Cell($i, $j).LeftPadding = 1.5
    }
}

# --- 2. Diary text for the finished registry work ---------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)

# Prefix the existing sentence with "1. " as its own run (force a run
# boundary by nudging a character property on just the inserted text,
# then resetting it back to its original/default value).
$start = $lastPara.Range.Start
$insertionPoint = $d.Range($start, $start)
$insertionPoint.InsertBefore("1. ")

$prefixRange = $d.Range($start, $start + 3)
$prefixRange.Bold = 1
$prefixRange.Bold = 0

# Append the new "2. ..." paragraph right after it.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = "2. Finished basic registry ajax post, path redirection, password comparision, authority code check and error dialog. "

# --- 3. Carry along the newly minted ListLabel character styles -------
for ($n = 127; $n -le 180; $n++) {
    $style = $d.Styles.Add("ListLabel$n", 2)
    $style.NameLocal = "ListLabel $n"
    $style.QuickStyle = $true
    $style.Font.NameBi = "OpenSymbol"
}
